$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.291.03"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.854.24"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "'" + "1.002"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'" + "314.56"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'" + "0.4602"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'" + "0.3710"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'" + "0.07296"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Value = "'" + "20.13"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "'" + "0.07832"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'" + "5.392"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.793.12"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").Value = "'" + "6.525"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "'" + "91.45"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'" + "1.003"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'" + "0.000008935"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'" + "14.75"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "27.311.06"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'" + "5.114"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "'" + "10.54"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "2.055.60"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "'" + "1.919"
$ws.Range("E25").Value = "  +4.14%  "
$ws.Range("D26").Value = "'" + "151.91"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'" + "18.47"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'" + "2.060"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'" + "116.01"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "'" + "0.08824"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'" + "3.093"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'" + "0.7728"
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("D34").Value = "'" + "1.170"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "'" + "4.515"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "'" + "2.764"
$ws.Range("E36").Value = "  +13.10%  "
$ws.Range("D37").Value = "'" + "1.080"
$ws.Range("D38").Value = "'" + "0.01954"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "'" + "0.05265"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "'" + "2.951"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "'" + "7.046"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "'" + "0.5125"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").Value = "'" + "0.1637"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "'" + "8.403"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "'" + "0.4793"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "'" + "10.42"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'" + "1.002"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'" + "102.40"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").Value = "'" + "0.06217"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'" + "65.71"
$ws.Range("E51").Value = "  +1.44%  "
